# Analysed experiments from 22-10-05 and 22-10-06
# Add "area um2", "uW/um2" and "W/m2" columns (D:F) computing the power
# density (in uW/um2 and W/m2) from the measured power (C) and the
# illuminated area (D, constant 1366 um2 for every row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New headers
$ws.Range("D1").Value = "area um2"
$ws.Range("E1").Value = "uW/um2"
$ws.Range("F1").Value = "W/m2"

# Area (constant) for every data row
$ws.Range("D2:D10").Value = 1366

# uW/um2 = power measured / area
$ws.Range("E2").Formula = "=C2/D2"
$ws.Range("E3:E10").Formula = "=C3/D3"

# W/m2 = uW/um2 * 10^6
$ws.Range("F2").Formula = "=E2*10^6"
$ws.Range("F3:F10").Formula = "=E3*10^6"

# Move/restore the active selection to F3, matching the authored workbook
$ws.Range("F3").Select()
